$d = $word.ActiveDocument

# Locate the "INTRODUCTION" heading paragraph robustly via Find, rather than
# a hard-coded paragraph index.
$rng = $d.Content
[void]$rng.Find.Execute("INTRODUCTION", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$introEnd = $rng.End

$targetIdx = -1
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($introEnd -ge $p.Range.Start) -and ($introEnd -le $p.Range.End)) {
        $targetIdx = $i
        break
    }
}

# The new paragraph is inserted right after the first empty paragraph that
# follows "INTRODUCTION" (i.e. immediately before the second, trailing empty
# paragraph that sits just before the section properties).
$insertBeforeIdx = $targetIdx + 2
$insertBeforePara = $d.Paragraphs.Item($insertBeforeIdx)
$insertRng = $insertBeforePara.Range
$insertRng.InsertParagraphBefore()

# Re-fetch: the freshly created empty paragraph now lives at $insertBeforeIdx,
# pushing the old trailing empty paragraph one slot further down.
$newPara = $d.Paragraphs.Item($insertBeforeIdx)
$newRng = $newPara.Range

# Build the paragraph exactly as authored (run boundaries + proofErr markers
# match the source Word session's grammar/spell-check squiggles).
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I have uploaded the Project Tweets .csv into </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>both of the technologies</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> used in the course. Vis a vie </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hbase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SparkSQL</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. And I have made some comparisons in my head on the performance of each off them. I think that for the task is it more convenient to use Spark SQL than </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hbase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> because it is visually easier to comprehend and one </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>doesn</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> not have to be </w:t></w:r><w:r><w:t xml:space="preserve">switching between terminals to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>eeecute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the commands.</w:t></w:r></w:p>
'@
[void]$newRng.InsertXML($xml)
